$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = -7.9075
$ws.Range("D12").Value = -5.989899999999998
$ws.Range("D18").Value = -8.640099999999999
$ws.Range("D37").Value = -8.148900000000001
$ws.Range("D55").Value = -8.969100000000001
$ws.Range("D68").Value = -6.951599999999996
$ws.Range("D77").Value = -6.317199999999997
$ws.Range("D78").Value = -7.849900000000002
